# Apply updated Bill of Materials changes to the BOM worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: D2 (Nexperia TVS diode) replaced with U1 (STMicroelectronics TVS diode)
$ws.Range("A13").Value = "STMicroelectronics"
$ws.Range("B13").Value = "USBLC6-2P6"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "U1"
$ws.Range("E13").Value = "17V Clamp 5A (8/20µs) Ipp Tvs Diode Surface Mount SOT-666"

# Row 15: U4 linear regulator description now specifies output voltage (1.8V)
$ws.Range("E15").Value = "Linear Voltage Regulator IC Positive Fixed 1 Output 700mA 4-XDFN (1x1), 1.8V"

# Row 16: U3 linear regulator description now specifies output voltage (3.3V) and is
# its own unique description (previously shared text with row 15)
$ws.Range("E16").Value = "Linear Voltage Regulator IC Positive Fixed 1 Output 700mA 4-XDFN (1x1), 3.3V"

# Row 20: R20 resistor description replaced with the standard Panasonic description text
$ws.Range("E20").Value = "2 kOhms ±1% 0.05W, 1/20W Chip Resistor 0201 (0603 Metric) Thick Film"

# Row 27: C15 added to the reference designators of the 0.1uF capacitor group
$ws.Range("D27").Value = "C1, C6, C7, C13, C15"

# Update the active selection to match the author's final cursor position
$ws.Range("B13").Select()
